$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '59.308.01'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '2.603.21'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.06'
$ws.Range('E5').Value = '  +4.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.49'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.44'
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').Value = '3.059.77'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '59.254.00'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.62'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.602.99'
$ws.Range('E16').Value = '  +1.70%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '341.26'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.13'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.40'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.60'
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.24'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('D28').Value = '0.0₃0745'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.68'
$ws.Range('E30').Value = '  +7.10%  '
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.74'
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.52'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +0.27%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '37.17'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.46'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.834'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.828'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '275.22'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0955'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('D47').Value = '1.951.08'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0224'
$ws.Range('E48').Value = '  +1.65%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.52'
$ws.Range('E49').Value = '  +3.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.52'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.02'
$ws.Range('E51').Value = '  -0.83%  '
